$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the Sanborns password from "Intermex202111" to "Intermex202201"
# (rows 4, 11 and 21 - column E, cadena = Sanborns)
$ws.Range("E4").Value = "Intermex202201"
$ws.Range("E11").Value = "Intermex202201"
$ws.Range("E21").Value = "Intermex202201"

# Update active selection to E7
$ws.Range("E7").Select()
